# Update the lattice-multiplication exercise table.
#
# The document holds a single 5-row x 3-column table; every cell's run
# contains 5 lines of text separated by <w:br/> line breaks:
#   1) "NN x MM"      - the multiplication problem
#   2) "  D1    D2"   - the two digits of the second factor, spaced out
#   3) "  ----"       - a constant separator
#   4) "D|    |"      - first digit of the first factor + lattice row
#   5) "D|    |"      - second digit of the first factor + lattice row
#
# Every one of the 15 cells gets new exercise content (new problem,
# same layout/formatting). We rewrite each cell's Range.Text wholesale,
# using Chr(11) (vertical tab) as the in-cell line-break character so
# Word re-creates the <w:br/> separated runs exactly as before.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

function Set-Cell($row, $col, $line1, $line2, $line3, $line4, $line5) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $line1 + $nl + $line2 + $nl + $line3 + $nl + $line4 + $nl + $line5
}

Set-Cell 1 1 "57 x 85" "  8    5" "  ----" "5|    |" "7|    |"
Set-Cell 1 2 "66 x 71" "  7    1" "  ----" "6|    |" "6|    |"
Set-Cell 1 3 "45 x 78" "  7    8" "  ----" "4|    |" "5|    |"

Set-Cell 2 1 "99 x 75" "  7    5" "  ----" "9|    |" "9|    |"
Set-Cell 2 2 "78 x 28" "  2    8" "  ----" "7|    |" "8|    |"
Set-Cell 2 3 "33 x 72" "  7    2" "  ----" "3|    |" "3|    |"

Set-Cell 3 1 "21 x 38" "  3    8" "  ----" "2|    |" "1|    |"
Set-Cell 3 2 "88 x 12" "  1    2" "  ----" "8|    |" "8|    |"
Set-Cell 3 3 "43 x 24" "  2    4" "  ----" "4|    |" "3|    |"

Set-Cell 4 1 "73 x 24" "  2    4" "  ----" "7|    |" "3|    |"
Set-Cell 4 2 "28 x 18" "  1    8" "  ----" "2|    |" "8|    |"
Set-Cell 4 3 "37 x 31" "  3    1" "  ----" "3|    |" "7|    |"

Set-Cell 5 1 "31 x 22" "  2    2" "  ----" "3|    |" "1|    |"
Set-Cell 5 2 "25 x 83" "  8    3" "  ----" "2|    |" "5|    |"
Set-Cell 5 3 "94 x 87" "  8    7" "  ----" "9|    |" "4|    |"

Write-Output "lattice table updated"
